$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(367, 44441, 0, 11, 177.9071648067281),
    @(368, 44442, 0, 10, 161.7337861879347),
    @(369, 44443, 0, 8, 129.3870289503477),
    @(370, 44444, 1, 7, 113.2136503315543),
    @(371, 44445, 3, 4, 64.69351447517387),
    @(372, 44446, 2, 6, 97.0402717127608),
    @(373, 44447, 0, 6, 97.0402717127608),
    @(374, 44448, 1, 7, 113.2136503315543)
)

foreach ($row in $data) {
    $r = $row[0]
    # Copy the format of the last existing data row (366) into the new row
    # so the date column keeps its date style/number format.
    $ws.Range("A366:D366").Copy($ws.Range("A$r`:D$r"))

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
